$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Mobile column (D) formatted as text so the leading-zero phone
# number is stored as a string rather than being coerced to a number.
$ws.Range("D3:D10").NumberFormat = "@"

# --- Append the submitted query rows (3-10) below the existing template
#     row (row 2 is the blank placeholder row already in the sheet). ---
$data = @(
    @("Gopal","Goyal","goyal11.gopal@gmail.com","09462976187","wassup"),
    @("Gopal","Goyal","goyal11.gopal@gmail.com","09462976187","wassup"),
    @("Gopal","Goyal","goyal11.gopal@gmail.com","09462976187","wassup"),
    @("Gopal","Goyal","goyal11.gopal@gmail.com","09462976187","wassup"),
    @("Gopal","Goyal","goyal11.gopal@gmail.com","09462976187","wassup"),
    $null,
    @("Gopal","Goyal","goyal11.gopal@gmail.com","09462976187","helloooo"),
    @("Gopal","Goyal","","09462976187","")
)

$r = 3
foreach ($row in $data) {
    if ($null -eq $row) {
        # Blank separator row (matches row 2's empty-string cells): write a
        # placeholder then clear it so the cell stays present-but-empty.
        foreach ($c in 1..5) {
            $cell = $ws.Cells.Item($r, $c)
            $cell.Value = "ZZZTEMPZZZ"
            $cell.Replace("ZZZTEMPZZZ", "") | Out-Null
        }
    } else {
        $c = 1
        foreach ($val in $row) {
            $ws.Cells.Item($r, $c).Value = $val
            $c++
        }
    }
    $r++
}

# --- Strip the bespoke alignment style that used to be applied to every
#     cell/column in the sheet; everything goes back to the default style.
#     (Column F only ever holds data in row 1, so it is handled on its
#     own to avoid materialising empty cells in F3:F10.) ---
$ws.Range("A1:E10").Style = "Normal"
$ws.Range("F1").Style = "Normal"
$ws.Columns("A:E").ClearFormats()
